$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.584.07"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "1.789.52"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.84"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.557"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.52"
$ws.Range("E8").Value = "  +5.60%  "
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0668"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").Value = "2.046.15"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.02"
$ws.Range("E13").Value = "  +10.13%  "
$ws.Range("D14").Value = "1.792.20"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "34.625.19"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.632"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.87"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.88"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "0.0₃0765"
$ws.Range("E20").Value = "  +3.30%  "
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.23"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.64"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.39"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.09"
$ws.Range("E27").Value = "  +1.73%  "
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("D35").Value = "1.443.51"
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0188"
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.626"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "83.15"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.82"
$ws.Range("E40").Value = "  +3.83%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.899"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.09"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0503"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.90"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("D47").Value = "1.945.28"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.96"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.99"
$ws.Range("E50").Value = "  +4.97%  "
$ws.Range("E51").Value = "  +4.60%  "
